$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 417
$ws.Range("F4").Value = 152
$ws.Range("F6").Value = 3793
$ws.Range("F7").Value = 222
$ws.Range("F8").Value = 2530
$ws.Range("F9").Value = 67
$ws.Range("F10").Value = 3061
$ws.Range("F11").Value = 1393
$ws.Range("F12").Value = 530
$ws.Range("F14").Value = 48
$ws.Range("F15").Value = 114
$ws.Range("F17").Value = 435
$ws.Range("F19").Value = 197
$ws.Range("F20").Value = 341
$ws.Range("F21").Value = 302
$ws.Range("F22").Value = 343
$ws.Range("F24").Value = 1387
$ws.Range("F25").Value = 40
$ws.Range("F26").Value = 1291
$ws.Range("F27").Value = 123
$ws.Range("F30").Value = 34
$ws.Range("F31").Value = 4213
$ws.Range("F32").Value = 3856
$ws.Range("F34").Value = 1
$ws.Range("F38").Value = 459
$ws.Range("F41").Value = 153
$ws.Range("F43").Value = 92
$ws.Range("F44").Value = 34
$ws.Range("F45").Value = 59

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 22
$ws.Range("F15").Value = 202

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 1023
$ws.Range("F4").Value = 2250

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 230
$ws.Range("F3").Value = 1023
$ws.Range("F5").Value = 417
$ws.Range("F8").Value = 152
$ws.Range("F10").Value = 3793
$ws.Range("F11").Value = 222
$ws.Range("F12").Value = 2530
$ws.Range("F13").Value = 67
$ws.Range("F14").Value = 3061
$ws.Range("F15").Value = 530
$ws.Range("F17").Value = 48
$ws.Range("F18").Value = 114
$ws.Range("F20").Value = 435
$ws.Range("F22").Value = 341
$ws.Range("F23").Value = 302
$ws.Range("F24").Value = 343
$ws.Range("F26").Value = 1387
$ws.Range("F27").Value = 40
$ws.Range("F28").Value = 1291
$ws.Range("F30").Value = 142
$ws.Range("F32").Value = 34
$ws.Range("F33").Value = 22
$ws.Range("F35").Value = 4213
$ws.Range("F36").Value = 3856
$ws.Range("F40").Value = 459
$ws.Range("F45").Value = 153
$ws.Range("F46").Value = 92
$ws.Range("F49").Value = 202

